# "handleDelete lesson is done"
#
# The "User Changes duration option / User Changes Start Time / User changes
# end time" Use-Case blocks used to be separated from the block above them
# (row 86) by a blank spacer row (row 87, style-only, no value). That spacer
# row (and the implicit blank row after it) is removed, so every block from
# row 90 downward shifts up by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unwanted blank rows 87:88 - this shifts rows 89:103 up by 2
# (89 was already blank/unused, 90->88, 91->89, ... 103->101), matching the
# new dimension A2:A101.
$ws.Rows("87:88").Delete()

# Restore the on-screen view: scrolled so row 57 is at the top, with D87
# selected as the active cell.
$ws.Range("D87").Select()
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 1
